$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.564.88"
$ws.Range("E2").Value = "  -2.10%  "

$ws.Range("D3").Value = "2.645.64"
$ws.Range("E3").Value = "  -3.34%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.27"
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("D9").Value = "2.645.79"
$ws.Range("E9").Value = "  -3.28%  "

$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.365"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.75%  "

$ws.Range("D15").Value = "3.118.87"
$ws.Range("E15").Value = "  -3.60%  "

$ws.Range("E16").Value = "  -3.90%  "

$ws.Range("D17").Value = "67.531.23"
$ws.Range("E17").Value = "  -2.14%  "

$ws.Range("D18").Value = "2.647.14"
$ws.Range("E18").Value = "  -2.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.90%  "

$ws.Range("E25").Value = "  -5.14%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.46%  "

$ws.Range("D28").Value = "2.777.36"
$ws.Range("E28").Value = "  -3.31%  "

$ws.Range("E29").Value = "  -3.63%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "556.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.26%  "

$ws.Range("E34").Value = "  -2.80%  "

$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "157.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.43%  "

$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("E41").Value = "  -4.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.93"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.36%  "

$ws.Range("D47").Value = "0.0₆0301"
$ws.Range("E47").Value = "  -3.04%  "

$ws.Range("E48").Value = "  -1.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.02%  "
